# Daily attendance processing - 2026-02-18 11:48:34 UTC
#
# 1. Column G ("Recorded By") is narrowed from width 31 -> 13 characters.
# 2. Every populated "Recorded By" cell (which held a recorder's name such
#    as "Miss Dina Nasr" or "Miss Dina Nasr, Administrator") is replaced
#    with the academic year "2025/2026".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Narrow column G (Recorded By) from 31 to 13 characters ------------
# Excel's ColumnWidth property (character units of the Normal style font)
# differs from the raw XML "width" attribute by a fixed per-font offset.
# Empirically, on this workbook's font, XML width = ColumnWidth + 0.83, so
# to land on an XML width of 13 we request ColumnWidth 12.17.
$ws.Columns.Item(7).ColumnWidth = 12.17

# --- 2. Replace recorder names in column G with "2025/2026" ---------------
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().StartsWith("Miss Dina Nasr")) {
        $cell.Value2 = "2025/2026"
    }
}
